# Trade #105 closed at 2026-02-17 15:58:57 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# --- Summary sheet -----------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 105     # Total Trades
$summary.Range("B9").Value = 38.1    # Win Rate %

# --- Strategy Status sheet ----------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 105      # MarketMaking Trades
$status.Range("G4").Value = 38.1     # MarketMaking Win Rate %

# --- New trade row data --------------------------------------------------
# Date/time are stored as plain text in this log (not Excel date serials),
# so they are entered with a leading apostrophe to force text and avoid
# Excel's automatic date/time parsing.
$tradeNum   = 105
$date       = "'2026-02-17"
$time       = "'15:58:51"
$strategy   = "MarketMaking"
$side       = "DOWN"
$entryPrice = 0.97
$exitPrice  = 0.97
$status_    = "CLOSED"
$plPct      = 0
$plDollar   = 0
$capAfter   = 100
$entrySlip  = 0
$exitSlip   = 0
$confidence = 0.6
$entryReason = "Normal spread capture: 19600 bps"
$exitReason  = "early_exit"
$duration    = 0.11

# --- All Trades sheet: append row 106 -----------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$row = 106
$allTrades.Cells.Item($row, 1).Value = $tradeNum
$allTrades.Cells.Item($row, 2).Value = $date
$allTrades.Cells.Item($row, 3).Value = $time
$allTrades.Cells.Item($row, 4).Value = $strategy
$allTrades.Cells.Item($row, 5).Value = $side
$allTrades.Cells.Item($row, 6).Value = $entryPrice
$allTrades.Cells.Item($row, 7).Value = $exitPrice
$allTrades.Cells.Item($row, 8).Value = $status_
$allTrades.Cells.Item($row, 9).Value = $plPct
$allTrades.Cells.Item($row, 10).Value = $plDollar
$allTrades.Cells.Item($row, 11).Value = $capAfter
$allTrades.Cells.Item($row, 12).Value = $entrySlip
$allTrades.Cells.Item($row, 13).Value = $exitSlip
$allTrades.Cells.Item($row, 14).Value = $confidence
$allTrades.Cells.Item($row, 15).Value = $entryReason
$allTrades.Cells.Item($row, 16).Value = $exitReason
$allTrades.Cells.Item($row, 17).Value = $duration

# --- MarketMaking sheet: append row 106 ---------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Cells.Item($row, 1).Value = $tradeNum
$mm.Cells.Item($row, 2).Value = $date
$mm.Cells.Item($row, 3).Value = $time
$mm.Cells.Item($row, 4).Value = $strategy
$mm.Cells.Item($row, 5).Value = $side
$mm.Cells.Item($row, 6).Value = $entryPrice
$mm.Cells.Item($row, 7).Value = $exitPrice
$mm.Cells.Item($row, 8).Value = $status_
$mm.Cells.Item($row, 9).Value = $plPct
$mm.Cells.Item($row, 10).Value = $plDollar
$mm.Cells.Item($row, 11).Value = $capAfter
$mm.Cells.Item($row, 12).Value = $entrySlip
$mm.Cells.Item($row, 13).Value = $exitSlip
$mm.Cells.Item($row, 14).Value = $confidence
$mm.Cells.Item($row, 15).Value = $entryReason
$mm.Cells.Item($row, 16).Value = $exitReason
$mm.Cells.Item($row, 17).Value = $duration
